# Auto-generated edit script applying the meteocat daily summary refresh
# (2026-02-07 07:49 run): refreshed extraction timestamps plus the small
# measurement deltas the diff captured for that pass.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text-valued cells (timestamps, hPa/°C/mm/km-h readings, etc.) ---
$ws.Range("E2").Value = "2026-02-07 07:47:36"
$ws.Range("O2").Value = "-1.9 °C"
$ws.Range("E3").Value = "2026-02-07 07:47:39"
$ws.Range("K3").Value = "0.0 MJ/m2"
$ws.Range("E4").Value = "2026-02-07 07:47:41"
$ws.Range("J4").Value = "1001.6 hPa"
$ws.Range("O4").Value = "11.0 °C"
$ws.Range("E5").Value = "2026-02-07 07:47:44"
$ws.Range("J5").Value = "1001.6 hPa"
$ws.Range("N5").Value = "5.7 °C 7:01 TU"
$ws.Range("O5").Value = "8.2 °C"
$ws.Range("E6").Value = "2026-02-07 07:47:46"
$ws.Range("J6").Value = "1003.2 hPa"
$ws.Range("K6").Value = "0.0 MJ/m2"
$ws.Range("E7").Value = "2026-02-07 07:47:48"
$ws.Range("J7").Value = "1002.9 hPa"
$ws.Range("E8").Value = "2026-02-07 07:47:51"
$ws.Range("K8").Value = "0.0 MJ/m2"
$ws.Range("E9").Value = "2026-02-07 07:47:53"
$ws.Range("O9").Value = "1.2 °C"
$ws.Range("E10").Value = "2026-02-07 07:47:55"
$ws.Range("E11").Value = "2026-02-07 07:47:58"
$ws.Range("J11").Value = "1005.8 hPa"
$ws.Range("O11").Value = "1.3 °C"
$ws.Range("E12").Value = "2026-02-07 07:48:00"
$ws.Range("O12").Value = "9.8 °C"
$ws.Range("E13").Value = "2026-02-07 07:48:03"
$ws.Range("O13").Value = "8.0 °C"
$ws.Range("E14").Value = "2026-02-07 07:48:05"
$ws.Range("I14").Value = "0.2 mm"
$ws.Range("N14").Value = "-8.3 °C 7:14 TU"
$ws.Range("O14").Value = "-6.0 °C"
$ws.Range("E15").Value = "2026-02-07 07:48:07"
$ws.Range("J15").Value = "1002.0 hPa"
$ws.Range("O15").Value = "5.6 °C"
$ws.Range("E16").Value = "2026-02-07 07:48:10"
$ws.Range("O16").Value = "2.7 °C"
$ws.Range("E17").Value = "2026-02-07 07:48:12"
$ws.Range("J17").Value = "1005.2 hPa"
$ws.Range("N17").Value = "2.2 °C 7:00 TU"
$ws.Range("O17").Value = "3.1 °C"
$ws.Range("E18").Value = "2026-02-07 07:48:15"
$ws.Range("E19").Value = "2026-02-07 07:48:17"
$ws.Range("J19").Value = "1006.6 hPa"
$ws.Range("O19").Value = "3.7 °C"
$ws.Range("E20").Value = "2026-02-07 07:48:19"
$ws.Range("K20").Value = "0.0 MJ/m2"
$ws.Range("N20").Value = "-5.8 °C 7:21 TU"
$ws.Range("E21").Value = "2026-02-07 07:48:21"
$ws.Range("J21").Value = "1002.4 hPa"
$ws.Range("O21").Value = "5.6 °C"
$ws.Range("E22").Value = "2026-02-07 07:48:24"
$ws.Range("K22").Value = "0.0 MJ/m2"
$ws.Range("L22").Value = "20.5 km/h - 336º 7:08 TU"
$ws.Range("M22").Value = "10.2 °C 7:09 TU"
$ws.Range("O22").Value = "6.7 °C"
$ws.Range("E23").Value = "2026-02-07 07:48:26"
$ws.Range("J23").Value = "1001.7 hPa"
$ws.Range("E24").Value = "2026-02-07 07:48:29"
$ws.Range("K24").Value = "0.0 MJ/m2"
$ws.Range("L24").Value = "73.1 km/h - 339º 7:14 TU"
$ws.Range("E25").Value = "2026-02-07 07:48:31"
$ws.Range("J25").Value = "1005.5 hPa"
$ws.Range("N25").Value = "0.2 °C 7:21 TU"
$ws.Range("E26").Value = "2026-02-07 07:48:33"
$ws.Range("O26").Value = "-2.6 °C"
$ws.Range("E27").Value = "2026-02-07 07:48:36"
$ws.Range("J27").Value = "1001.5 hPa"
$ws.Range("L27").Value = "31.0 km/h - 0º 7:20 TU"
$ws.Range("M27").Value = "11.6 °C 7:01 TU"
$ws.Range("O27").Value = "9.0 °C"
$ws.Range("E28").Value = "2026-02-07 07:48:38"
$ws.Range("J28").Value = "1004.4 hPa"
$ws.Range("N28").Value = "1.3 °C 7:29 TU"
$ws.Range("O28").Value = "2.9 °C"
$ws.Range("E29").Value = "2026-02-07 07:48:41"
$ws.Range("K29").Value = "0.0 MJ/m2"
$ws.Range("O29").Value = "10.6 °C"
$ws.Range("E30").Value = "2026-02-07 07:48:43"
$ws.Range("K30").Value = "0.1 MJ/m2"
$ws.Range("M30").Value = "-2.8 °C 7:19 TU"
$ws.Range("O30").Value = "-4.8 °C"
$ws.Range("E31").Value = "2026-02-07 07:48:45"
$ws.Range("J31").Value = "1006.1 hPa"
$ws.Range("E32").Value = "2026-02-07 07:48:47"
$ws.Range("J32").Value = "1004.7 hPa"
$ws.Range("E33").Value = "2026-02-07 07:48:50"
$ws.Range("O33").Value = "6.7 °C"
$ws.Range("E34").Value = "2026-02-07 07:48:52"
$ws.Range("K34").Value = "0.0 MJ/m2"
$ws.Range("E35").Value = "2026-02-07 07:48:56"
$ws.Range("O35").Value = "-6.2 °C"
$ws.Range("E36").Value = "2026-02-07 07:48:58"
$ws.Range("J36").Value = "1006.9 hPa"
$ws.Range("K36").Value = "0.0 MJ/m2"
$ws.Range("N36").Value = "2.6 °C 7:13 TU"
$ws.Range("O36").Value = "4.4 °C"

# --- Percentage-looking cells: force Text format first so the COM layer
#     keeps them as literal strings ("68%") instead of converting them to
#     a numeric 0.68 with a Percent number format. ---
$percentCells = @("H12", "H13", "H20", "H21", "H22", "H23", "H24", "H27", "H30")
foreach ($pc in $percentCells) {
    $ws.Range($pc).NumberFormat = "@"
}
$ws.Range("H12").Value = "68%"
$ws.Range("H13").Value = "85%"
$ws.Range("H20").Value = "81%"
$ws.Range("H21").Value = "79%"
$ws.Range("H22").Value = "88%"
$ws.Range("H23").Value = "97%"
$ws.Range("H24").Value = "80%"
$ws.Range("H27").Value = "92%"
$ws.Range("H30").Value = "79%"
